$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("approvalDetails")

# Update the approver remarks text for the existing last row (row 17):
# "Grievance approval required" -> "Funds needed to resolve the issue"
$ws.Range("E17").Value = "Funds needed to resolve the issue"

# Append the new row of flow data (row 18)
$ws.Range("A18").Value = "juniorAssistant"
$ws.Range("B18").Value = "ACCOUNTS"
$ws.Range("C18").Value = "Junior Assistant"
$ws.Range("D18").Value = "D.Subramanyam"
$ws.Range("E18").Value = "Sanctioned and shall grievance be processed"

# Widen column E so the new, longer remarks text fits
$ws.Columns.Item(5).ColumnWidth = 40.5

# Reset the view: scroll back so column A is visible again, and move the
# active selection down to B22 (matches the saved view state)
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B22").Select()
